# rill-analysis: Control page display.
#
# - Point the REST endpoint (the "_settings" sheet's B1 "url" row) at the
#   new rill-analysis-web service instead of the old standalone saiku one,
#   and drop the now-unused "[Time].[2011]" / "时间" helper labels that used
#   to sit next to it in D1/E1.
# - Widen column B on "_settings" so the longer URL still fits.
# - Show the "_settings" sheet (not "_input") when the workbook opens.

$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("_settings")
$input = $wb.Worksheets.Item("_input")

# Swap the saiku REST endpoint for the rill-analysis-web one.
$settings.Range("B1").Value = "http://10.81.21.140:8280/rill-analysis-web/rest/"

# These two helper cells ("时间" / "[Time].[2011]") are no longer needed.
$settings.Range("D1").Value = ""
$settings.Range("E1").Value = ""

# The new URL is longer than the old one - widen column B (from 40.5 to
# ~53.875 characters) so it still fits.
$settings.Columns.Item(2).ColumnWidth = 53.142857142857146

# Display "_settings" as the active tab instead of "_input".
$settings.Activate()
